# Update "Línea 141" horarios workbook: new scrape adds one record
# (23_HERNANDEZ, Hora_Scrap 16:34:35, Hora_Llegada 17:02, 28 min, LP1912)
# to the LP1912 sheet's data table, which is kept sorted by Hora_Llegada.
# All three sheets get their "Última actualización" timestamp refreshed.

$wb = $excel.ActiveWorkbook

# ---- Sheet "LP1912": header + full (re-sorted) data table ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 16:34:35"
$ws1.Range("A3").Value = "Total filas: 134"

$data = @(
    @("12:41:18","12:41","15X38_ABASTO",0,"LP1912"),
    @("12:58:39","13:01","10_OLMOS",3,"LP1912"),
    @("12:38:18","13:04","23_HERNANDEZ",26,"LP1912"),
    @("12:58:39","13:21","16_SANTA ANA",23,"LP1912"),
    @("12:58:39","13:33","23_HERNANDEZ",35,"LP1912"),
    @("13:35:25","13:36","17_ROMERO",1,"LP1912"),
    @("12:46:01","13:38","23_HERNANDEZ",52,"LP1912"),
    @("12:46:01","13:43","14_ABASTO",57,"LP1912"),
    @("13:35:25","13:51","16_SANTA ANA",16,"LP1912"),
    @("13:54:15","13:55","17_ROMERO",1,"LP1912"),
    @("13:35:25","14:01","16_SANTA ANA",26,"LP1912"),
    @("13:35:25","14:03","23_HERNANDEZ",28,"LP1912"),
    @("13:54:15","14:05","14_ABASTO",11,"LP1912"),
    @("12:58:39","14:11","15_ABASTO",73,"LP1912"),
    @("13:54:15","14:13","16_SANTA ANA",19,"LP1912"),
    @("13:35:25","14:14","10_OLMOS",39,"LP1912"),
    @("12:58:39","14:16","27_EL RETIRO",78,"LP1912"),
    @("12:58:39","14:20","215C_EL PATO",82,"LP1912"),
    @("12:38:18","14:21","215C_EL PATO",103,"LP1912"),
    @("12:41:18","14:25","11_ETCHEVERRY",104,"LP1912"),
    @("13:54:15","14:31","15_ABASTO",37,"LP1912"),
    @("14:23:38","14:33","23_HERNANDEZ",10,"LP1912"),
    @("12:58:39","14:33","215C_LA PLATA",95,"LP1912"),
    @("12:46:01","14:34","215C_LA PLATA",108,"LP1912"),
    @("13:35:25","14:34","10_OLMOS",59,"LP1912"),
    @("12:41:18","14:37","16_P MOR-SANTA ANA",116,"LP1912"),
    @("13:35:25","14:38","23_HERNANDEZ",63,"LP1912"),
    @("12:41:18","14:40","17X38_ROMERO",119,"LP1912"),
    @("12:58:39","14:44","215B_EL PATO",106,"LP1912"),
    @("12:46:01","14:45","215B_EL PATO",119,"LP1912"),
    @("14:23:38","14:49","16_SANTA ANA",26,"LP1912"),
    @("14:43:48","14:51","15_ABASTO",8,"LP1912"),
    @("12:58:39","14:53","14_ABASTO",115,"LP1912"),
    @("12:58:39","14:53","215A_LA PLATA",115,"LP1912"),
    @("14:43:48","14:54","10_OLMOS",11,"LP1912"),
    @("13:35:25","14:56","215A_LA PLATA",81,"LP1912"),
    @("14:57:55","14:58","16_P MOR-SANTA ANA",1,"LP1912"),
    @("14:43:48","15:00","16_SANTA ANA",17,"LP1912"),
    @("13:35:25","15:01","81_EL PELIGRO",86,"LP1912"),
    @("14:57:55","15:02","16_SANTA ANA",5,"LP1912"),
    @("13:54:15","15:02","215A_LA PLATA",68,"LP1912"),
    @("14:23:38","15:03","23_HERNANDEZ",40,"LP1912"),
    @("13:35:25","15:04","14_ABASTO",89,"LP1912"),
    @("14:23:38","15:04","16_SANTA ANA",41,"LP1912"),
    @("14:43:48","15:04","215A_LA PLATA",21,"LP1912"),
    @("13:54:15","15:05","14_ABASTO",71,"LP1912"),
    @("13:35:25","15:17","16_P MOR-SANTA ANA",102,"LP1912"),
    @("14:23:38","15:24","215C_EL PATO",61,"LP1912"),
    @("13:35:25","15:24","11_ETCHEVERRY",109,"LP1912"),
    @("13:54:15","15:25","11_ETCHEVERRY",91,"LP1912"),
    @("13:54:15","15:25","16_P MOR-SANTA ANA",91,"LP1912"),
    @("13:35:25","15:25","215C_EL PATO",110,"LP1912"),
    @("14:57:55","15:33","16_SANTA ANA",36,"LP1912"),
    @("13:54:15","15:36","17X38_ROMERO",102,"LP1912"),
    @("14:23:38","15:44","14_ABASTO",81,"LP1912"),
    @("14:23:38","15:45","215C_LA PLATA",82,"LP1912"),
    @("15:40:34","15:48","215C_LA PLATA",8,"LP1912"),
    @("14:57:55","15:51","215C_LA PLATA",54,"LP1912"),
    @("14:43:48","15:53","215C_LA PLATA",70,"LP1912"),
    @("14:23:38","15:55","27_EL RETIRO",92,"LP1912"),
    @("14:43:48","15:56","27_EL RETIRO",73,"LP1912"),
    @("15:40:34","15:57","16_SANTA ANA",17,"LP1912"),
    @("15:59:48","16:00","27_EL RETIRO",1,"LP1912"),
    @("14:23:38","16:01","15_ABASTO",98,"LP1912"),
    @("15:40:34","16:03","23_HERNANDEZ",23,"LP1912"),
    @("15:59:48","16:04","14_ABASTO",5,"LP1912"),
    @("15:40:34","16:14","10_OLMOS",34,"LP1912"),
    @("14:57:55","16:16","10_OLMOS",79,"LP1912"),
    @("14:23:38","16:20","10_OLMOS",117,"LP1912"),
    @("15:40:34","16:21","16_SANTA ANA",41,"LP1912"),
    @("14:57:55","16:22","15_ABASTO",85,"LP1912"),
    @("16:22:52","16:22","11_ETCHEVERRY",0,"LP1912"),
    @("14:43:48","16:24","11_ETCHEVERRY",101,"LP1912"),
    @("14:43:48","16:24","215_ALUAR",101,"LP1912"),
    @("14:43:48","16:25","215B_LP-P MOR-1 Y 57",102,"LP1912"),
    @("15:40:34","16:25","215_ALUAR",45,"LP1912"),
    @("15:40:34","16:26","215B_LP-P MOR-1 Y 57",46,"LP1912"),
    @("14:43:48","16:31","16_P MOR-SANTA ANA",108,"LP1912"),
    @("15:40:34","16:32","23_HERNANDEZ",52,"LP1912"),
    @("16:32:08","16:33","16_P MOR-SANTA ANA",1,"LP1912"),
    @("14:57:55","16:33","225_GOMEZ",96,"LP1912"),
    @("14:43:48","16:33","23_HERNANDEZ",110,"LP1912"),
    @("14:43:48","16:34","225_GOMEZ",111,"LP1912"),
    @("14:43:48","16:36","17X38_ROMERO",113,"LP1912"),
    @("15:40:34","16:38","15_ABASTO",58,"LP1912"),
    @("14:57:55","16:41","17X38_ROMERO",104,"LP1912"),
    @("16:16:23","16:44","23_HERNANDEZ",28,"LP1912"),
    @("16:22:52","16:50","23_HERNANDEZ",28,"LP1912"),
    @("15:59:48","16:51","16_SANTA ANA",52,"LP1912"),
    @("14:57:55","16:53","10_OLMOS",116,"LP1912"),
    @("16:16:23","16:54","10_OLMOS",38,"LP1912"),
    @("16:28:39","16:55","23_HERNANDEZ",27,"LP1912"),
    @("14:57:55","16:56","215C_EL PATO",119,"LP1912"),
    @("15:40:34","16:57","215C_LA PLATA",77,"LP1912"),
    @("15:40:34","16:57","215C_EL PATO",77,"LP1912"),
    @("16:32:08","17:00","23_HERNANDEZ",28,"LP1912"),
    @("16:22:52","17:01","16_SANTA ANA",39,"LP1912"),
    @("16:34:35","17:02","23_HERNANDEZ",28,"LP1912"),
    @("15:40:34","17:03","14_ABASTO",83,"LP1912"),
    @("15:59:48","17:03","215C_LA PLATA",64,"LP1912"),
    @("15:59:48","17:04","14_ABASTO",65,"LP1912"),
    @("16:32:08","17:05","14_ABASTO",33,"LP1912"),
    @("16:16:23","17:07","15_ABASTO",51,"LP1912"),
    @("16:33:53","17:13","10_OLMOS",40,"LP1912"),
    @("16:22:52","17:14","10_OLMOS",52,"LP1912"),
    @("15:40:34","17:17","17_ROMERO",97,"LP1912"),
    @("15:40:34","17:24","11_ETCHEVERRY",104,"LP1912"),
    @("16:32:08","17:25","11_ETCHEVERRY",53,"LP1912"),
    @("15:59:48","17:27","15_ABASTO",88,"LP1912"),
    @("16:16:23","17:32","23_HERNANDEZ",76,"LP1912"),
    @("16:22:52","17:33","23_HERNANDEZ",71,"LP1912"),
    @("15:40:34","17:34","10_OLMOS",114,"LP1912"),
    @("16:32:08","17:35","23_HERNANDEZ",63,"LP1912"),
    @("15:40:34","17:35","16_P MOR-SANTA ANA",115,"LP1912"),
    @("15:40:34","17:36","27_EL RETIRO",116,"LP1912"),
    @("15:40:34","17:38","17X38_ROMERO",118,"LP1912"),
    @("15:59:48","17:40","27_EL RETIRO",101,"LP1912"),
    @("15:59:48","17:44","215B_EL PATO",105,"LP1912"),
    @("16:32:08","17:45","215B_EL PATO",73,"LP1912"),
    @("16:22:52","17:45","17X38_ROMERO",83,"LP1912"),
    @("16:16:23","17:46","17X38_ROMERO",90,"LP1912"),
    @("16:22:52","17:46","27_EL RETIRO",84,"LP1912"),
    @("15:59:48","17:47","17X38_ROMERO",108,"LP1912"),
    @("16:33:53","17:47","27_EL RETIRO",74,"LP1912"),
    @("16:16:23","17:48","27_EL RETIRO",92,"LP1912"),
    @("16:28:39","17:49","27_EL RETIRO",81,"LP1912"),
    @("15:59:48","17:50","215_EL PELIGRO",111,"LP1912"),
    @("16:16:23","17:51","215_EL PELIGRO",95,"LP1912"),
    @("16:22:52","18:02","17_ROMERO",100,"LP1912"),
    @("16:16:23","18:03","17_ROMERO",107,"LP1912"),
    @("16:16:23","18:04","14_ABASTO",108,"LP1912"),
    @("16:32:08","18:05","14_ABASTO",93,"LP1912"),
    @("16:28:39","18:24","11_ETCHEVERRY",116,"LP1912"),
    @("16:32:08","18:25","11_ETCHEVERRY",113,"LP1912")
)

$r = 6
foreach ($row in $data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---- Sheet "LP1912-215": timestamp only ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 16:34:35"

# ---- Sheet "6203-6173": timestamp only ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 16:34:35"

Write-Output "done"
